# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2310"
#   "<header>_new" -> "<header>_FV2404"
# and turn the sheet's used range into a real Excel Table with the header
# row frozen, so the renamed headers stay visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row: columns A-J carried the "_old" suffix, column K
#    is the constant "diff" column, and columns L-U carried the "_new"
#    suffix. Replace those suffixes with the respective format-version tag.
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2310"
}
$ws.Cells.Item(1, 11).Value = "diff"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2404"
}

# 2) Freeze the header row so it stays visible when scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Convert the used range A1:U61 into a native Excel Table ("Table1")
#    with an auto filter on the header row.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U61"), $null, 1)
$lo.Name = "Table1"

$ws.Range("A1").Select()
